# Atualização de bases das ligas, do dia: 22-05-2024 às 20:16
#
# This script re-applies a re-shuffle of several match rows (their B:AB
# data, i.e. everything except the rank number in column A and the
# Div/Date columns C/D which are shared/unchanged) plus a couple of
# plain odds-value corrections on two not-yet-played fixtures.
#
# Because several of the edits are *cyclic* row permutations (row X's
# data becomes row Y's, row Y's becomes row Z's, row Z's becomes row
# X's, etc.), every source row's B:AB values are snapshotted into a
# variable BEFORE any writes happen, and only then written out to their
# destination rows. This avoids clobbering data that is still needed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Get-RowData($row) {
    return $ws.Range("B$row`:AB$row").Value2
}

function Set-RowData($row, $data) {
    $ws.Range("B$row`:AB$row").Value2 = $data
}

# ---------------------------------------------------------------------
# Group 1: rows 208-212 (cyclic permutation)
#   208 <- 211, 209 <- 210, 210 <- 212, 211 <- 208, 212 <- 209
# ---------------------------------------------------------------------
$row208 = Get-RowData 208
$row209 = Get-RowData 209
$row210 = Get-RowData 210
$row211 = Get-RowData 211
$row212 = Get-RowData 212

Set-RowData 208 $row211
Set-RowData 209 $row210
Set-RowData 210 $row212
Set-RowData 211 $row208
Set-RowData 212 $row209

# ---------------------------------------------------------------------
# Group 2: rows 424-426 (cyclic permutation)
#   424 <- 425, 425 <- 426, 426 <- 424
# ---------------------------------------------------------------------
$row424 = Get-RowData 424
$row425 = Get-RowData 425
$row426 = Get-RowData 426

Set-RowData 424 $row425
Set-RowData 425 $row426
Set-RowData 426 $row424

# ---------------------------------------------------------------------
# Group 3: rows 429-431 (cyclic permutation)
#   429 <- 430, 430 <- 431, 431 <- 429
# ---------------------------------------------------------------------
$row429 = Get-RowData 429
$row430 = Get-RowData 430
$row431 = Get-RowData 431

Set-RowData 429 $row430
Set-RowData 430 $row431
Set-RowData 431 $row429

# ---------------------------------------------------------------------
# Group 4: row 444 - plain odds-value corrections
# ---------------------------------------------------------------------
$ws.Cells.Item(444, 13).Value = 2.3      # M444 oddH
$ws.Cells.Item(444, 15).Value = 3.75     # O444 oddA
$ws.Cells.Item(444, 17).Value = 1.975    # Q444 oddAHH
$ws.Cells.Item(444, 18).Value = 1.875    # R444 oddAHA
$ws.Cells.Item(444, 20).Value = 1.85     # T444 oddAHOver
$ws.Cells.Item(444, 21).Value = 2        # U444 oddAHUnder

# ---------------------------------------------------------------------
# Group 5: row 445 - plain odds-value corrections
# ---------------------------------------------------------------------
$ws.Cells.Item(445, 13).Value = 2.5      # M445 oddH
$ws.Cells.Item(445, 15).Value = 3.25     # O445 oddA
$ws.Cells.Item(445, 16).Value = -0.25    # P445 Ah
$ws.Cells.Item(445, 17).Value = 2.125    # Q445 oddAHH
$ws.Cells.Item(445, 18).Value = 1.75     # R445 oddAHA
